$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 369.04166
$ws.Range("J17").Value = 367.69565
$ws.Range("L17").Value = 1103.08695
$ws.Range("N17").Value = -1439.08695

$ws.Range("H39").Value = 130.47058
$ws.Range("I39").Value = 66.07692
$ws.Range("J39").Value = 339.75
$ws.Range("K39").Value = 198.23076
$ws.Range("L39").Value = 1019.25
$ws.Range("M39").Value = 97.76924
$ws.Range("N39").Value = -1611.25

$ws.Range("H76").Value = 3971350
$ws.Range("I76").Value = 4276531
$ws.Range("J76").Value = 3999.5
$ws.Range("K76").Value = 4276531
$ws.Range("L76").Value = 3999.5
$ws.Range("M76").Value = -4276216
$ws.Range("N76").Value = -4629.5

$ws.Range("H79").Value = 3971350
$ws.Range("I79").Value = 4276531
$ws.Range("J79").Value = 3999.5
$ws.Range("K79").Value = 4276531
$ws.Range("L79").Value = 3999.5
$ws.Range("M79").Value = -4275439
$ws.Range("N79").Value = -6183.5

$ws.Range("H101").Value = 6988.2354
$ws.Range("I101").Value = 512.8889
$ws.Range("J101").Value = 14273
$ws.Range("K101").Value = 1538.6667
$ws.Range("L101").Value = 42819
$ws.Range("M101").Value = 83.33329999999978
$ws.Range("N101").Value = -46063

$ws.Range("H123").Value = 96661.73
$ws.Range("J123").Value = 96661.73
$ws.Range("L123").Value = 96661.73
$ws.Range("N123").Value = -106461.73

$ws.Range("H137").Value = 1581.1305
$ws.Range("I137").Value = 900.4286
$ws.Range("J137").Value = 1878.9375
$ws.Range("K137").Value = 2701.2858
$ws.Range("L137").Value = 5636.8125
$ws.Range("M137").Value = -151.2857999999997
$ws.Range("N137").Value = -10736.8125

$ws.Range("H138").Value = 6581287
$ws.Range("I138").Value = 2676.8635
$ws.Range("J138").Value = 9261461
$ws.Range("K138").Value = 8030.5905
$ws.Range("L138").Value = 27784383
$ws.Range("M138").Value = -2890.5905
$ws.Range("N138").Value = -27794663

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws.Range("H32").Value = 2512.0254
$ws.Range("I32").Value = 1364.2667
$ws.Range("J32").Value = 6136.5264
$ws.Range("K32").Value = 1364.2667
$ws.Range("L32").Value = 6136.5264
$ws.Range("M32").Value = -1077.2667
$ws.Range("N32").Value = -6710.5264

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H101").Value = 25000
$ws.Range("J101").Value = 25000
$ws.Range("L101").Value = 25000
$ws.Range("N101").Value = -31490

$ws.Range("H122").Value = 9388.308000000001
$ws.Range("I122").Value = 10747.454
$ws.Range("J122").Value = 1913
$ws.Range("K122").Value = 32242.362
$ws.Range("L122").Value = 5739
$ws.Range("M122").Value = -29792.362
$ws.Range("N122").Value = -10639

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H139").Value = 66614
$ws.Range("J139").Value = 66614
$ws.Range("L139").Value = 66614
$ws.Range("N139").Value = -76894

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 2661
$ws.Range("I10").Value = 2005
$ws.Range("J10").Value = 2825
$ws.Range("K10").Value = 2005
$ws.Range("L10").Value = 2825
$ws.Range("M10").Value = -1865
$ws.Range("N10").Value = -3105

$ws.Range("H134").Value = 2762.1282
$ws.Range("I134").Value = 1640.25
$ws.Range("K134").Value = 4920.75
$ws.Range("M134").Value = -2385.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1397.7972
$ws.Range("I31").Value = 1035.7567
$ws.Range("J31").Value = 1759.8379
$ws.Range("K31").Value = 1035.7567
$ws.Range("L31").Value = 1759.8379
$ws.Range("M31").Value = -740.7566999999999
$ws.Range("N31").Value = -2349.8379

$ws.Range("H34").Value = 1397.7972
$ws.Range("I34").Value = 1035.7567
$ws.Range("J34").Value = 1759.8379
$ws.Range("K34").Value = 1035.7567
$ws.Range("L34").Value = 1759.8379
$ws.Range("M34").Value = -833.7566999999999
$ws.Range("N34").Value = -2163.8379

$ws.Range("H62").Value = 24770.6
$ws.Range("I62").Value = 35950
$ws.Range("J62").Value = 8001.5
$ws.Range("K62").Value = 35950
$ws.Range("L62").Value = 8001.5
$ws.Range("M62").Value = -35326
$ws.Range("N62").Value = -9249.5

$ws.Range("H65").Value = 24770.6
$ws.Range("I65").Value = 35950
$ws.Range("J65").Value = 8001.5
$ws.Range("K65").Value = 179750
$ws.Range("L65").Value = 40007.5
$ws.Range("M65").Value = -176630
$ws.Range("N65").Value = -46247.5

$ws.Range("H96").Value = 14610.083
$ws.Range("J96").Value = 14610.083
$ws.Range("L96").Value = 14610.083
$ws.Range("N96").Value = -20102.083

$ws.Range("H99").Value = 3907310
$ws.Range("I99").Value = 4465211.5
$ws.Range("K99").Value = 4465211.5
$ws.Range("M99").Value = -4463713.5

$ws.Range("H106").Value = 50000
$ws.Range("J106").Value = 50000
$ws.Range("L106").Value = 50000
$ws.Range("N106").Value = -52524

$ws.Range("H126").Value = 3907310
$ws.Range("I126").Value = 4465211.5
$ws.Range("K126").Value = 13395634.5
$ws.Range("M126").Value = -13393164.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 862.28
$ws.Range("I68").Value = 736.6769399999999
$ws.Range("J68").Value = 1095.5428
$ws.Range("K68").Value = 2210.03082
$ws.Range("L68").Value = 3286.6284
$ws.Range("M68").Value = -1399.03082
$ws.Range("N68").Value = -4908.6284

$ws.Range("H71").Value = 862.28
$ws.Range("I71").Value = 736.6769399999999
$ws.Range("J71").Value = 1095.5428
$ws.Range("K71").Value = 6630.09246
$ws.Range("L71").Value = 9859.885199999999
$ws.Range("M71").Value = -2574.09246
$ws.Range("N71").Value = -17971.8852

$ws.Range("H97").Value = 613.7143
$ws.Range("I97").Value = 659.2
$ws.Range("J97").Value = 500
$ws.Range("K97").Value = 1977.6
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = -1481.6
$ws.Range("N97").Value = -2492

$ws.Range("H107").Value = 812.9091
$ws.Range("I107").Value = 181.80952
$ws.Range("J107").Value = 1202.7059
$ws.Range("K107").Value = 545.4285599999999
$ws.Range("L107").Value = 3608.1177
$ws.Range("M107").Value = 1374.57144
$ws.Range("N107").Value = -7448.1177

$ws.Range("H131").Value = 2348.189
$ws.Range("J131").Value = 2545.4197
$ws.Range("L131").Value = 7636.259099999999
$ws.Range("N131").Value = -17716.2591

$ws.Range("H137").Value = 3745468.2
$ws.Range("I137").Value = 7146970
$ws.Range("J137").Value = 82312.16
$ws.Range("K137").Value = 21440910
$ws.Range("L137").Value = 246936.48
$ws.Range("M137").Value = -21435810
$ws.Range("N137").Value = -257136.48

$ws.Range("H139").Value = 1945
$ws.Range("I139").Value = 1784.2106
$ws.Range("K139").Value = 5352.6318
$ws.Range("M139").Value = -212.6318000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H51").Value = 57000
$ws.Range("J51").Value = 57000
$ws.Range("L51").Value = 57000
$ws.Range("N51").Value = -58018

$ws.Range("H53").Value = 20727.273
$ws.Range("I53").Value = 8000
$ws.Range("K53").Value = 8000
$ws.Range("M53").Value = -7369

$ws.Range("H70").Value = 5604.077
$ws.Range("I70").Value = 6929.0713
$ws.Range("J70").Value = 4058.25
$ws.Range("K70").Value = 6929.0713
$ws.Range("L70").Value = 4058.25
$ws.Range("M70").Value = -6659.0713
$ws.Range("N70").Value = -4598.25

$ws.Range("H73").Value = 5604.077
$ws.Range("I73").Value = 6929.0713
$ws.Range("J73").Value = 4058.25
$ws.Range("K73").Value = 6929.0713
$ws.Range("L73").Value = 4058.25
$ws.Range("M73").Value = -5993.0713
$ws.Range("N73").Value = -5930.25

$ws.Range("H80").Value = 2550
$ws.Range("I80").Value = 2500
$ws.Range("J80").Value = 2600
$ws.Range("K80").Value = 2500
$ws.Range("L80").Value = 2600
$ws.Range("M80").Value = -1502
$ws.Range("N80").Value = -4596

$ws.Range("H83").Value = 2550
$ws.Range("I83").Value = 2500
$ws.Range("J83").Value = 2600
$ws.Range("K83").Value = 12500
$ws.Range("L83").Value = 13000
$ws.Range("M83").Value = -7508
$ws.Range("N83").Value = -22984

$ws.Range("H102").Value = 1753.4286
$ws.Range("I102").Value = 1213.9231
$ws.Range("J102").Value = 3312
$ws.Range("K102").Value = 1213.9231
$ws.Range("L102").Value = 3312
$ws.Range("M102").Value = 408.0769
$ws.Range("N102").Value = -6556

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 24935
$ws.Range("J104").Value = 24935
$ws.Range("L104").Value = 24935
$ws.Range("N104").Value = -31923

$ws.Range("H123").Value = 38214.5
$ws.Range("J123").Value = 38214.5
$ws.Range("L123").Value = 38214.5
$ws.Range("N123").Value = -48014.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 29867.334
$ws.Range("J101").Value = 29867.334
$ws.Range("L101").Value = 29867.334
$ws.Range("N101").Value = -36357.334

$ws.Range("H123").Value = 29857
$ws.Range("J123").Value = 29857
$ws.Range("L123").Value = 29857
$ws.Range("N123").Value = -39657

$ws.Range("H132").Value = 11629557
$ws.Range("I132").Value = 13159108
$ws.Range("J132").Value = 4968
$ws.Range("K132").Value = 39477324
$ws.Range("L132").Value = 14904
$ws.Range("M132").Value = -39474794
$ws.Range("N132").Value = -19964

$ws.Range("H136").Value = 14537494
$ws.Range("I136").Value = 18574588
$ws.Range("K136").Value = 55723764
$ws.Range("M136").Value = -55721214
